$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (keys) and column B (values) for rows 1-16 to the new
# simplified key/value naming scheme. Values that look like plain numbers
# or the word "true" are prefixed with a leading apostrophe so Excel keeps
# them stored as text (matching the shared-string cells in the target
# workbook) instead of coercing them to numeric/boolean cells; the style
# is reset to "Normal" right after so no stray cell formatting is left
# behind.

$ws.Range("A2").Value = "gris_value"
$ws.Range("B2").Value = "0,30"

$ws.Range("A3").Value = "gris_min"
$ws.Range("B3").Value = "isento"

$ws.Range("A4").Value = "toll_value"
$ws.Range("B4").Value = "4,25"

$ws.Range("A5").Value = "delivery_value"
$ws.Range("B5").Value = "isento"

$ws.Range("A6").Value = "tas_value"
$ws.Range("B6").Value = "isento"

$ws.Range("A7").Value = "other_fee_value"
$ws.Range("B7").Value = "isento"

$ws.Range("A8").Value = "other_fee_fraction"
$ws.Range("B8").Value = "isento"

$ws.Range("A9").Value = "cubic_factor"
$ws.Range("B9").Value = "'280"
$ws.Range("B9").Style = "Normal"

$ws.Range("A10").Value = "tda_value"
$ws.Range("B10").Value = "'50"
$ws.Range("B10").Style = "Normal"

$ws.Range("A11").Value = "tda_min"
$ws.Range("B11").Value = "351,85"

$ws.Range("A12").Value = "tda_max"
$ws.Range("B12").Value = "882,43"

$ws.Range("A13").Value = "trt_value"
$ws.Range("B13").Value = "isento"

$ws.Range("A14").Value = "trt_min"
$ws.Range("B14").Value = "isento"

$ws.Range("A15").Value = "icms_tax_enable"
$ws.Range("B15").Value = "'true"
$ws.Range("B15").Style = "Normal"

$ws.Range("A16").Value = "iss_tax_enable"
$ws.Range("B16").Value = "'true"
$ws.Range("B16").Style = "Normal"

# Rows 17 and 18 (old ICMS/ISS rows) are no longer needed; remove them
# (delete bottom-up so row numbers of earlier rows stay stable).
$ws.Rows("18").Delete()
$ws.Rows("17").Delete()
